# Update the "Förändrad" (column C) date for all data rows (2-83)
# from 2023-09-14 (serial 45183) to 2023-09-15 (serial 45184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C83").Value = 45184
